# Add data for 2022-10-15: update sheet name, header label, and October/Total rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab / workbook sheet name.
$ws.Name = "Through 2022-10-07"

# Update the "October (through 10-06)" label to "October (through 10-07)".
$ws.Range("A11").Value = "October (through 10-07)"

# Update the October row (row 11) values for the columns that changed.
$ws.Range("C11").Value = 10
$ws.Range("E11").Value = 16
$ws.Range("G11").Value = 35
$ws.Range("H11").Value = 51
$ws.Range("I11").Value = 22

# Update the Total row (row 12) values for the columns that changed.
$ws.Range("C12").Value = 439
$ws.Range("E12").Value = 564
$ws.Range("G12").Value = 936
$ws.Range("H12").Value = 1298
$ws.Range("I12").Value = 1303
